# "lots of work on port-level landings"
# Row 2 currently has the "MONTEREY AREA TOTALS" label sitting in column B
# (port) with nothing in column A (species). Move that label into A2 (so
# it lines up as a port-level row) and put the "Totals" label into B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "MONTEREY AREA TOTALS"
$ws.Range("B2").Value = "Totals"

# Column A no longer needs its own narrower best-fit width now that it
# holds the same kind of long label text as column B - widen it to match.
$ws.Columns("A:A").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# Selection moves from the single cell A2 to the whole of column A.
$ws.Columns("A:A").Select()
